$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the is_locked_lbl (D1) and is_enabled_lbl (E1) columns by shifting
# order_by (was F1) and rem (was G1) two columns to the left, then put the
# new tenant_id_lbl validation text where rem used to sit (now the last col).
$ws.Range("D1").Value = '<%=comment.order_by%>'
$ws.Range("E1").Value = '<%=comment.rem%>'
$ws.Range("F1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# Drop the now-vacated trailing column so the used range shrinks back to F1.
$ws.Range("G1").ClearContents()
